$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45047.33333333334, "FX_IDC:USDUAH", 36.5684, 36.5685, 36.5681, 36.5681, 0),
    @(45078.33333333334, "FX_IDC:USDUAH", 36.5681, 36.5681, 36.565,  36.565,  0),
    @(45110.33333333334, "FX_IDC:USDUAH", 36.565,  36.565,  36.565,  36.565,  0)
)

$startRow = 292
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    # Copy formatting (style) from the last existing data row (291) so the
    # new rows match the sheet's established look (style index reused).
    $ws.Range("A291:G291").Copy($ws.Range("A$row`:G$row"))

    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
}
